$d = $word.ActiveDocument

# This chapter's "Titre1" headings are manually numbered with roman
# numerals typed directly into the text ("II - ...", "III - ..."). Two
# of the headings (a second "II" and a "III") were mis-numbered and are
# being bumped up by one (II -> III, III -> IV) to restore the correct
# sequence I, II, III, IV.

function Bump-HeadingNumber($searchPhrase, $oldNumber, $newNumber) {
    $rng = $d.Content
    $found = $rng.Find.Execute($searchPhrase, $true, $true, $false, $false, `
                                $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Heading not found: $searchPhrase"
    }

    # $rng now covers just the matched phrase; narrow it down to the
    # leading roman-numeral run only, and swap that run's text.
    $numRange = $d.Range($rng.Start, $rng.Start + $oldNumber.Length)
    if ($numRange.Text -ne $oldNumber) {
        throw "Unexpected heading numeral: [$($numRange.Text)]"
    }
    $numRange.Text = $newNumber
}

Bump-HeadingNumber "II – Différence de deux nombres relatifs" "II" "III"
Bump-HeadingNumber "III – Simplification d'une suite de sommes" "III" "IV"
